# Auto-generated edit script applying the Bahamut_Profits leve-profit
# recalculation update across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H46").Value = 5463.6
$ws.Range("J46").Value = 5463.6
$ws.Range("L46").Value = 16390.8
$ws.Range("N46").Value = -16628.8

$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()

$ws.Range("H60").Value = 5463.6
$ws.Range("J60").Value = 5463.6
$ws.Range("L60").Value = 16390.8
$ws.Range("N60").Value = -17358.8

$ws.Range("H86").Value = 3194.1177
$ws.Range("I86").Value = 2378.5715
$ws.Range("J86").Value = 7000
$ws.Range("K86").Value = 2378.5715
$ws.Range("L86").Value = 7000
$ws.Range("M86").Value = -1255.5715
$ws.Range("N86").Value = -9246

$ws.Range("H89").Value = 3194.1177
$ws.Range("I89").Value = 2378.5715
$ws.Range("J89").Value = 7000
$ws.Range("K89").Value = 11892.8575
$ws.Range("L89").Value = 35000
$ws.Range("M89").Value = -6276.8575
$ws.Range("N89").Value = -46232

$ws.Range("H138").Value = 3066.55
$ws.Range("J138").Value = 3277.4092
$ws.Range("L138").Value = 9832.2276
$ws.Range("N138").Value = -20112.2276

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11362.506
$ws.Range("I32").Value = 6737.831
$ws.Range("K32").Value = 6737.831
$ws.Range("M32").Value = -6450.831

$ws.Range("H61").Value = 2273.1428
$ws.Range("I61").Value = 2273.1428
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2273.1428
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -2061.1428
$ws.Range("N61").ClearContents()

$ws.Range("H136").Value = 2273.1428
$ws.Range("I136").Value = 2273.1428
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 6819.428400000001
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -4269.428400000001
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 47996.08
$ws.Range("I20").Value = 101349.836
$ws.Range("J20").Value = 2264.2856
$ws.Range("K20").Value = 101349.836
$ws.Range("L20").Value = 2264.2856
$ws.Range("M20").Value = -101102.836
$ws.Range("N20").Value = -2758.2856

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2415.5745
$ws.Range("I31").Value = 2475.2104
$ws.Range("J31").Value = 2163.7778
$ws.Range("K31").Value = 2475.2104
$ws.Range("L31").Value = 2163.7778
$ws.Range("M31").Value = -2180.2104
$ws.Range("N31").Value = -2753.7778

$ws.Range("H33").Value = 2137
$ws.Range("I33").Value = 715.5
$ws.Range("K33").Value = 715.5
$ws.Range("M33").Value = -336.5

$ws.Range("H34").Value = 2415.5745
$ws.Range("I34").Value = 2475.2104
$ws.Range("J34").Value = 2163.7778
$ws.Range("K34").Value = 2475.2104
$ws.Range("L34").Value = 2163.7778
$ws.Range("M34").Value = -2273.2104
$ws.Range("N34").Value = -2567.7778

$ws.Range("H68").Value = 15632.777
$ws.Range("J68").Value = 15632.777
$ws.Range("L68").Value = 15632.777
$ws.Range("N68").Value = -17130.777

$ws.Range("H71").Value = 15632.777
$ws.Range("J71").Value = 15632.777
$ws.Range("L71").Value = 46898.331
$ws.Range("N71").Value = -54386.331

$ws.Range("H74").Value = 12320
$ws.Range("I74").Value = 5500
$ws.Range("J74").Value = 13684
$ws.Range("K74").Value = 5500
$ws.Range("L74").Value = 13684
$ws.Range("M74").Value = -4626
$ws.Range("N74").Value = -15432

$ws.Range("H77").Value = 12320
$ws.Range("I77").Value = 5500
$ws.Range("J77").Value = 13684
$ws.Range("K77").Value = 16500
$ws.Range("L77").Value = 41052
$ws.Range("M77").Value = -12132
$ws.Range("N77").Value = -49788

$ws.Range("H122").Value = 621.5
$ws.Range("I122").Value = 459.33334
$ws.Range("J122").Value = 718.8
$ws.Range("K122").Value = 1378.00002
$ws.Range("L122").Value = 2156.4
$ws.Range("M122").Value = 1071.99998
$ws.Range("N122").Value = -7056.4

$ws.Range("H134").Value = 15626062
$ws.Range("I134").Value = 926.76
$ws.Range("J134").Value = 71430110
$ws.Range("K134").Value = 2780.28
$ws.Range("L134").Value = 214290330
$ws.Range("M134").Value = -245.2799999999997
$ws.Range("N134").Value = -214295400

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 6812.5
$ws.Range("I80").Value = 3166.6667
$ws.Range("J80").Value = 9000
$ws.Range("K80").Value = 9500.000100000001
$ws.Range("L80").Value = 27000
$ws.Range("M80").Value = -8564.000100000001
$ws.Range("N80").Value = -28872

$ws.Range("H83").Value = 6812.5
$ws.Range("I83").Value = 3166.6667
$ws.Range("J83").Value = 9000
$ws.Range("K83").Value = 28500.0003
$ws.Range("L83").Value = 81000
$ws.Range("M83").Value = -23820.0003
$ws.Range("N83").Value = -90360

$ws.Range("H113").Value = 18081.648
$ws.Range("I113").Value = 732.36365
$ws.Range("J113").Value = 22230.39
$ws.Range("K113").Value = 2197.09095
$ws.Range("L113").Value = 66691.17
$ws.Range("M113").Value = -27.09094999999979
$ws.Range("N113").Value = -71031.17

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 732474.25
$ws.Range("I122").Value = 1463047.5
$ws.Range("J122").Value = 1901
$ws.Range("K122").Value = 4389142.5
$ws.Range("L122").Value = 5703
$ws.Range("M122").Value = -4386692.5
$ws.Range("N122").Value = -10603

$ws.Range("H128").Value = 33000
$ws.Range("J128").Value = 33000
$ws.Range("L128").Value = 33000
$ws.Range("N128").Value = -42960

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 981.25
$ws.Range("I22").Value = 400
$ws.Range("J22").Value = 1064.2858
$ws.Range("K22").Value = 400
$ws.Range("L22").Value = 1064.2858
$ws.Range("M22").Value = -105
$ws.Range("N22").Value = -1654.2858

$ws.Range("H27").Value = 981.25
$ws.Range("I27").Value = 400
$ws.Range("J27").Value = 1064.2858
$ws.Range("K27").Value = 400
$ws.Range("L27").Value = 1064.2858
$ws.Range("M27").Value = -293
$ws.Range("N27").Value = -1278.2858

$ws.Range("H40").Value = 1444898.6
$ws.Range("I40").Value = 2021798
$ws.Range("J40").Value = 2650
$ws.Range("K40").Value = 2021798
$ws.Range("L40").Value = 2650
$ws.Range("M40").Value = -2021662
$ws.Range("N40").Value = -2922

$ws.Range("H122").Value = 1627.1111
$ws.Range("I122").Value = 1406.7693
$ws.Range("J122").Value = 2200
$ws.Range("K122").Value = 4220.3079
$ws.Range("L122").Value = 6600
$ws.Range("M122").Value = -1770.3079
$ws.Range("N122").Value = -11500

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1739.4828
$ws.Range("I81").Value = 1259.2858
$ws.Range("K81").Value = 2518.5716
$ws.Range("M81").Value = -1457.5716

$ws.Range("H84").Value = 1739.4828
$ws.Range("I84").Value = 1259.2858
$ws.Range("K84").Value = 12592.858
$ws.Range("M84").Value = -7288.858

$ws.Range("H109").Value = 15433.333
$ws.Range("J109").Value = 15433.333
$ws.Range("L109").Value = 15433.333
$ws.Range("N109").Value = -18207.333

